$wb = $excel.ActiveWorkbook

# --- 1. Rename first sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Area Minimized Results"

# --- 2. Update the "Leakage Power" label text (nW -> W) on every block ---
$ws1.Range("B6").Value = "Leakage Power (W)"
$ws1.Range("B13").Value = "Leakage Power (W)"
$ws1.Range("B20").Value = "Leakage Power (W)"
$ws1.Range("B27").Value = "Leakage Power (W)"
$ws1.Range("B34").Value = "Leakage Power (W)"
$ws1.Range("B41").Value = "Leakage Power (W)"

# --- 3. Row 6 (D0 block) gets an actual leakage-power measurement ---
$rng6 = $ws1.Range("C6:M6")
$rng6.HorizontalAlignment = -4108
$rng6.NumberFormat = "0.00E+00"
$rng6.Merge()
$ws1.Range("C6").Value = 0.000000000895

# --- 4. Remaining blocks (D1,D2,D3,D4,D6) get blank merged placeholder cells ---
$rng13 = $ws1.Range("C13:M13")
$rng13.HorizontalAlignment = -4108
$rng13.Merge()

$rng20 = $ws1.Range("C20:M20")
$rng20.HorizontalAlignment = -4108
$rng20.Merge()

$rng27 = $ws1.Range("C27:M27")
$rng27.HorizontalAlignment = -4108
$rng27.Merge()

$rng34 = $ws1.Range("C34:M34")
$rng34.HorizontalAlignment = -4108
$rng34.Merge()

$rng41 = $ws1.Range("C41:M41")
$rng41.HorizontalAlignment = -4108
$rng41.Merge()

# --- 5. "Area Plot" sheet: add a ratio column (D) ---
$ws2.Range("D11").Formula = "=B11/C11"
$ws2.Range("D12").Formula = "=B12/C12"
$ws2.Range("D13").Formula = "=B13/C13"

# Column A width on the plot sheet
$ws2.Columns.Item(1).ColumnWidth = 13.6

# --- 6. Chart: give the (hidden) value axis a rotated title ---
$co = $ws2.ChartObjects().Item(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.HasTitle = $true
$valAx.AxisTitle.Text = "Area (a.u.)"

# --- 7. Selection / active tab bookkeeping ---
$ws2.Range("D14").Select()
$ws1.Select()
$ws1.Range("C14:M14").Select()

$wb.Save()
